$wb = $excel.ActiveWorkbook

# Row 3 on ALC (diff hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 56266.668  # H3: 0 -> 56266.668
$ws.Cells.Item(3, 10).Value = 56266.668  # J3: 0 -> 56266.668
$ws.Cells.Item(3, 12).Value = 56266.668  # L3: 0 -> 56266.668
$ws.Cells.Item(3, 14).Value = -56494.668  # N3: None -> -56494.668

# Row 6 on ALC (diff hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 10000164  # H6: 25000250 -> 10000164
$ws.Cells.Item(6, 9).Value = 204.75  # I6: 0 -> 204.75
$ws.Cells.Item(6, 10).Value = 50000000  # J6: 25000250 -> 50000000
$ws.Cells.Item(6, 11).Value = 614.25  # K6: 0 -> 614.25
$ws.Cells.Item(6, 12).Value = 150000000  # L6: 75000750 -> 150000000
$ws.Cells.Item(6, 13).Value = -502.25  # M6: None -> -502.25
$ws.Cells.Item(6, 14).Value = -150000224  # N6: -75000974 -> -150000224

# Row 62 on ALC (diff hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 400001730  # H62: 333335200 -> 400001730
$ws.Cells.Item(62, 10).Value = 0  # J62: 2500 -> 0
$ws.Cells.Item(62, 12).Value = 0  # L62: 2500 -> 0
$ws.Cells.Item(62, 14).ClearContents()  # N62: -3748 -> (removed)

# Row 65 on ALC (diff hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 400001730  # H65: 333335200 -> 400001730
$ws.Cells.Item(65, 10).Value = 0  # J65: 2500 -> 0
$ws.Cells.Item(65, 12).Value = 0  # L65: 12500 -> 0
$ws.Cells.Item(65, 14).ClearContents()  # N65: -18740 -> (removed)

# Row 102 on ALC (diff hunk 4)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(102, 8).Value = 56266.668  # H102: 0 -> 56266.668
$ws.Cells.Item(102, 10).Value = 56266.668  # J102: 0 -> 56266.668
$ws.Cells.Item(102, 12).Value = 56266.668  # L102: 0 -> 56266.668
$ws.Cells.Item(102, 14).Value = -62756.668  # N102: None -> -62756.668

# Row 133 on ALC (diff hunk 5)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value = 56433.332  # H133: 59650 -> 56433.332
$ws.Cells.Item(133, 10).Value = 56433.332  # J133: 59650 -> 56433.332
$ws.Cells.Item(133, 12).Value = 56433.332  # L133: 59650 -> 56433.332
$ws.Cells.Item(133, 14).Value = -66553.33199999999  # N133: -69770 -> -66553.33199999999

# Row 132 on ARM (diff hunk 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 20855726  # H132: 21762556 -> 20855726
$ws.Cells.Item(132, 9).Value = 29412976  # I132: 30304328 -> 29412976
$ws.Cells.Item(132, 10).Value = 73833.42999999999  # J132: 79596.46000000001 -> 73833.42999999999
$ws.Cells.Item(132, 11).Value = 88238928  # K132: 90912984 -> 88238928
$ws.Cells.Item(132, 12).Value = 221500.29  # L132: 238789.38 -> 221500.29
$ws.Cells.Item(132, 13).Value = -88236398  # M132: -90910454 -> -88236398
$ws.Cells.Item(132, 14).Value = -226560.29  # N132: -243849.38 -> -226560.29

# Row 133 on ARM (diff hunk 7)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(133, 8).Value = 0  # H133: 30000 -> 0
$ws.Cells.Item(133, 10).Value = 0  # J133: 30000 -> 0
$ws.Cells.Item(133, 12).Value = 0  # L133: 30000 -> 0
$ws.Cells.Item(133, 14).ClearContents()  # N133: -35060 -> (removed)

# Row 132 on BSM (diff hunk 8)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(132, 8).Value = 51260  # H132: 54390 -> 51260
$ws.Cells.Item(132, 10).Value = 51260  # J132: 54390 -> 51260
$ws.Cells.Item(132, 12).Value = 51260  # L132: 54390 -> 51260
$ws.Cells.Item(132, 14).Value = -61380  # N132: -64510 -> -61380

# Row 134 on BSM (diff hunk 9)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 6773.5854  # H134: 6970.846 -> 6773.5854
$ws.Cells.Item(134, 9).Value = 3470.1936  # I134: 3442.3438 -> 3470.1936
$ws.Cells.Item(134, 10).Value = 17014.1  # J134: 23101.143 -> 17014.1
$ws.Cells.Item(134, 11).Value = 10410.5808  # K134: 10327.0314 -> 10410.5808
$ws.Cells.Item(134, 12).Value = 51042.3  # L134: 69303.429 -> 51042.3
$ws.Cells.Item(134, 13).Value = -7875.5808  # M134: -7792.0314 -> -7875.5808
$ws.Cells.Item(134, 14).Value = -56112.3  # N134: -74373.429 -> -56112.3

# Row 16 on CRP (diff hunk 10)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2184.4443  # H16: 1690.8334 -> 2184.4443
$ws.Cells.Item(16, 9).Value = 1690  # I16: 1448 -> 1690
$ws.Cells.Item(16, 10).Value = 2580  # J16: 1864.2858 -> 2580
$ws.Cells.Item(16, 11).Value = 1690  # K16: 1448 -> 1690
$ws.Cells.Item(16, 12).Value = 2580  # L16: 1864.2858 -> 2580
$ws.Cells.Item(16, 13).Value = -1403  # M16: -1161 -> -1403
$ws.Cells.Item(16, 14).Value = -3154  # N16: -2438.2858 -> -3154

# Row 20 on CRP (diff hunk 11)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 39882.25  # H20: 36640.777 -> 39882.25
$ws.Cells.Item(20, 9).Value = 0  # I20: 10709 -> 0
$ws.Cells.Item(20, 11).Value = 0  # K20: 10709 -> 0
$ws.Cells.Item(20, 13).ClearContents()  # M20: -10473 -> (removed)

# Row 30 on CRP (diff hunk 12)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(30, 8).Value = 39882.25  # H30: 36640.777 -> 39882.25
$ws.Cells.Item(30, 9).Value = 0  # I30: 10709 -> 0
$ws.Cells.Item(30, 11).Value = 0  # K30: 10709 -> 0
$ws.Cells.Item(30, 13).ClearContents()  # M30: -10618 -> (removed)

# Row 31 on CRP (diff hunk 13)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2901769.8  # H31: 2501778 -> 2901769.8
$ws.Cells.Item(31, 9).Value = 3403153.8  # I31: 3573276.5 -> 3403153.8
$ws.Cells.Item(31, 10).Value = 269503.5  # J31: 120669.78 -> 269503.5
$ws.Cells.Item(31, 11).Value = 3403153.8  # K31: 3573276.5 -> 3403153.8
$ws.Cells.Item(31, 12).Value = 269503.5  # L31: 120669.78 -> 269503.5
$ws.Cells.Item(31, 13).Value = -3402858.8  # M31: -3572981.5 -> -3402858.8
$ws.Cells.Item(31, 14).Value = -270093.5  # N31: -121259.78 -> -270093.5

# Row 34 on CRP (diff hunk 14)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2901769.8  # H34: 2501778 -> 2901769.8
$ws.Cells.Item(34, 9).Value = 3403153.8  # I34: 3573276.5 -> 3403153.8
$ws.Cells.Item(34, 10).Value = 269503.5  # J34: 120669.78 -> 269503.5
$ws.Cells.Item(34, 11).Value = 3403153.8  # K34: 3573276.5 -> 3403153.8
$ws.Cells.Item(34, 12).Value = 269503.5  # L34: 120669.78 -> 269503.5
$ws.Cells.Item(34, 13).Value = -3402951.8  # M34: -3573074.5 -> -3402951.8
$ws.Cells.Item(34, 14).Value = -269907.5  # N34: -121073.78 -> -269907.5

# Row 62 on CRP (diff hunk 15)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 0  # H62: 3101.2 -> 0
$ws.Cells.Item(62, 9).Value = 0  # I62: 3000 -> 0
$ws.Cells.Item(62, 10).Value = 0  # J62: 3253 -> 0
$ws.Cells.Item(62, 11).Value = 0  # K62: 3000 -> 0
$ws.Cells.Item(62, 12).Value = 0  # L62: 3253 -> 0
$ws.Cells.Item(62, 13).ClearContents()  # M62: -2376 -> (removed)
$ws.Cells.Item(62, 14).ClearContents()  # N62: -4501 -> (removed)

# Row 65 on CRP (diff hunk 16)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 0  # H65: 3101.2 -> 0
$ws.Cells.Item(65, 9).Value = 0  # I65: 3000 -> 0
$ws.Cells.Item(65, 10).Value = 0  # J65: 3253 -> 0
$ws.Cells.Item(65, 11).Value = 0  # K65: 15000 -> 0
$ws.Cells.Item(65, 12).Value = 0  # L65: 16265 -> 0
$ws.Cells.Item(65, 13).ClearContents()  # M65: -11880 -> (removed)
$ws.Cells.Item(65, 14).ClearContents()  # N65: -22505 -> (removed)

# Row 86 on CRP (diff hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 21009.066  # H86: 25067.334 -> 21009.066
$ws.Cells.Item(86, 9).Value = 21170  # I86: 26250 -> 21170
$ws.Cells.Item(86, 10).Value = 20767.666  # J86: 24121.2 -> 20767.666
$ws.Cells.Item(86, 11).Value = 21170  # K86: 26250 -> 21170
$ws.Cells.Item(86, 12).Value = 20767.666  # L86: 24121.2 -> 20767.666
$ws.Cells.Item(86, 13).Value = -20047  # M86: -25127 -> -20047
$ws.Cells.Item(86, 14).Value = -23013.666  # N86: -26367.2 -> -23013.666

# Row 89 on CRP (diff hunk 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 21009.066  # H89: 25067.334 -> 21009.066
$ws.Cells.Item(89, 9).Value = 21170  # I89: 26250 -> 21170
$ws.Cells.Item(89, 10).Value = 20767.666  # J89: 24121.2 -> 20767.666
$ws.Cells.Item(89, 11).Value = 105850  # K89: 131250 -> 105850
$ws.Cells.Item(89, 12).Value = 103838.33  # L89: 120606 -> 103838.33
$ws.Cells.Item(89, 13).Value = -100234  # M89: -125634 -> -100234
$ws.Cells.Item(89, 14).Value = -115070.33  # N89: -131838 -> -115070.33

# Row 113 on CRP (diff hunk 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 2184.4443  # H113: 1690.8334 -> 2184.4443
$ws.Cells.Item(113, 9).Value = 1690  # I113: 1448 -> 1690
$ws.Cells.Item(113, 10).Value = 2580  # J113: 1864.2858 -> 2580
$ws.Cells.Item(113, 11).Value = 1690  # K113: 1448 -> 1690
$ws.Cells.Item(113, 12).Value = 2580  # L113: 1864.2858 -> 2580
$ws.Cells.Item(113, 13).Value = 480  # M113: 722 -> 480
$ws.Cells.Item(113, 14).Value = -6920  # N113: -6204.2858 -> -6920

# Row 128 on CRP (diff hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(128, 8).Value = 39882.25  # H128: 36640.777 -> 39882.25
$ws.Cells.Item(128, 9).Value = 0  # I128: 10709 -> 0
$ws.Cells.Item(128, 11).Value = 0  # K128: 10709 -> 0
$ws.Cells.Item(128, 13).ClearContents()  # M128: -5729 -> (removed)

# Row 132 on CRP (diff hunk 21)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 6063194  # H132: 4653236 -> 6063194
$ws.Cells.Item(132, 9).Value = 2015.6538  # I132: 1906.9 -> 2015.6538
$ws.Cells.Item(132, 10).Value = 28576142  # J132: 15387072 -> 28576142
$ws.Cells.Item(132, 11).Value = 6046.9614  # K132: 5720.700000000001 -> 6046.9614
$ws.Cells.Item(132, 12).Value = 85728426  # L132: 46161216 -> 85728426
$ws.Cells.Item(132, 13).Value = -3516.9614  # M132: -3190.700000000001 -> -3516.9614
$ws.Cells.Item(132, 14).Value = -85733486  # N132: -46166276 -> -85733486

# Row 7 on CUL (diff hunk 22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 82.52631  # H7: 75 -> 82.52631
$ws.Cells.Item(7, 9).Value = 70  # I7: 50 -> 70
$ws.Cells.Item(7, 10).Value = 84  # J7: 100 -> 84
$ws.Cells.Item(7, 11).Value = 210  # K7: 150 -> 210
$ws.Cells.Item(7, 12).Value = 252  # L7: 300 -> 252
$ws.Cells.Item(7, 13).Value = -98  # M7: -38 -> -98
$ws.Cells.Item(7, 14).Value = -476  # N7: -524 -> -476

# Row 126 on GSM (diff hunk 23)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 55589560  # H126: 15162583 -> 55589560
$ws.Cells.Item(126, 9).Value = 51005  # I126: 13973.125 -> 51005
$ws.Cells.Item(126, 10).Value = 166666670  # J126: 55558876 -> 166666670
$ws.Cells.Item(126, 11).Value = 153015  # K126: 41919.375 -> 153015
$ws.Cells.Item(126, 12).Value = 500000010  # L126: 166676628 -> 500000010
$ws.Cells.Item(126, 13).Value = -150545  # M126: -39449.375 -> -150545
$ws.Cells.Item(126, 14).Value = -500004950  # N126: -166681568 -> -500004950

# Row 132 on GSM (diff hunk 24)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 74206.42999999999  # H132: 61145.707 -> 74206.42999999999
$ws.Cells.Item(132, 9).Value = 2456.125  # I132: 1968.1818 -> 2456.125
$ws.Cells.Item(132, 10).Value = 169873.5  # J132: 169637.83 -> 169873.5
$ws.Cells.Item(132, 11).Value = 7368.375  # K132: 5904.5454 -> 7368.375
$ws.Cells.Item(132, 12).Value = 509620.5  # L132: 508913.49 -> 509620.5
$ws.Cells.Item(132, 13).Value = -4838.375  # M132: -3374.5454 -> -4838.375
$ws.Cells.Item(132, 14).Value = -514680.5  # N132: -513973.49 -> -514680.5

# Row 133 on GSM (diff hunk 25)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(133, 8).Value = 57950  # H133: 48000 -> 57950
$ws.Cells.Item(133, 10).Value = 57950  # J133: 48000 -> 57950
$ws.Cells.Item(133, 12).Value = 57950  # L133: 48000 -> 57950
$ws.Cells.Item(133, 14).Value = -68070  # N133: -58120 -> -68070

# Row 132 on LTW (diff hunk 26)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3790305  # H132: 10420252 -> 3790305
$ws.Cells.Item(132, 9).Value = 5557381.5  # I132: 41668120 -> 5557381.5
$ws.Cells.Item(132, 10).Value = 3713.2856  # J132: 4296.3335 -> 3713.2856
$ws.Cells.Item(132, 11).Value = 16672144.5  # K132: 125004360 -> 16672144.5
$ws.Cells.Item(132, 12).Value = 11139.8568  # L132: 12889.0005 -> 11139.8568
$ws.Cells.Item(132, 13).Value = -16669614.5  # M132: -125001830 -> -16669614.5
$ws.Cells.Item(132, 14).Value = -16199.8568  # N132: -17949.0005 -> -16199.8568

# Row 132 on WVR (diff hunk 27)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 54499530  # H132: 210813170 -> 54499530
$ws.Cells.Item(132, 9).Value = 77586760  # I132: 321429340 -> 77586760
$ws.Cells.Item(132, 10).Value = 6675979  # J132: 17234886 -> 6675979
$ws.Cells.Item(132, 11).Value = 232760280  # K132: 964288020 -> 232760280
$ws.Cells.Item(132, 12).Value = 20027937  # L132: 51704658 -> 20027937
$ws.Cells.Item(132, 13).Value = -232757750  # M132: -964285490 -> -232757750
$ws.Cells.Item(132, 14).Value = -20032997  # N132: -51709718 -> -20032997
